# QA_Assignment import rework:
# the "Open URL" row (row 3) should link to https://keystoneglobalnetwork.com/,
# not the stray https://www.amazon.com/ value it had, and the leftover
# hyperlink-only cell at F18 goes away.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the old hyperlink that lived on F18 (leaving the plain cell behind).
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$F$18') {
        $h.Delete()
    }
}
$ws.Range("F18").Value = ""

# Row 3 ("Open URL") should carry the real URL + working hyperlink instead.
$ws.Range("F3").Value = "https://keystoneglobalnetwork.com/"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://keystoneglobalnetwork.com/")

# Leave the cursor where it was left in the saved file.
$ws.Range("E18").Select() | Out-Null
